# Update the "Key Stage 4 (KS4) destinations" and "Key Stage 5 (KS5) destinations"
# rows with refreshed National Pupil Database permalinks and new period ranges.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 11: Key Stage 4 (KS4) destinations
$ws.Range("B11").Value = "<a href = 'https://explore-education-statistics.service.gov.uk/data-tables/permalink/86d98c08-f1f2-469f-5a0f-08dcf4c92e7d'>National Pupil Database</a>"
$ws.Range("C11").Value = "Aug 2022 -  Jul 2023 (21/22 learners) (24/10/24)"
$ws.Range("D11").Value = "Aug 2023 -  Jul 2024 (22/23 learners) (Oct 25)"

# Row 12: Key Stage 5 (KS5) destinations
$ws.Range("B12").Value = "<a href = 'https://explore-education-statistics.service.gov.uk/data-tables/permalink/3c47beb5-97bc-4cb2-5a12-08dcf4c92e7d'>National Pupil Database</a>"
$ws.Range("C12").Value = "Aug 2022 -  Jul 2023 (21/22 learners) (24/10/24)"
$ws.Range("D12").Value = "Aug 2023 -  Jul 2024 (22/23 learners) (Oct 25)"

# Update selected cell to match the saved view state.
$ws.Range("B12").Select()

$wb.Save()
